$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '37.903.42'
Set-TextValue $ws.Range("E2") '  +1.81%  '
Set-TextValue $ws.Range("D3") '2.052.68'
Set-TextValue $ws.Range("E3") '  +1.25%  '
Set-TextValue $ws.Range("E4") '  +0.12%  '
Set-TextValue $ws.Range("D5") '229.70'
Set-TextValue $ws.Range("E5") '  +1.04%  '
Set-TextValue $ws.Range("E6") '  +2.17%  '
Set-TextValue $ws.Range("E7") '  +5.23%  '
Set-TextValue $ws.Range("E8") '  +0.01%  '
Set-TextValue $ws.Range("E9") '  +1.96%  '
Set-TextValue $ws.Range("E10") '  +2.38%  '
Set-TextValue $ws.Range("D11") '0.103'
Set-TextValue $ws.Range("E11") '  +0.58%  '
Set-TextValue $ws.Range("D12") '2.355.96'
Set-TextValue $ws.Range("E12") '  +1.26%  '
Set-TextValue $ws.Range("E13") '  +2.39%  '
Set-TextValue $ws.Range("D14") '20.72'
Set-TextValue $ws.Range("E14") '  +1.57%  '
Set-TextValue $ws.Range("E15") '  +2.35%  '
Set-TextValue $ws.Range("D16") '0.750'
Set-TextValue $ws.Range("E16") '  +1.02%  '
Set-TextValue $ws.Range("D17") '2.050.73'
Set-TextValue $ws.Range("E17") '  +1.16%  '
Set-TextValue $ws.Range("D18") '37.839.21'
Set-TextValue $ws.Range("E19") '  -2.56%  '
Set-TextValue $ws.Range("E20") '  +0.62%  '
Set-TextValue $ws.Range("D21") '0.0₃0832'
Set-TextValue $ws.Range("E21") '  +1.27%  '
Set-TextValue $ws.Range("D22") '224.70'
Set-TextValue $ws.Range("E22") '  +0.05%  '
Set-TextValue $ws.Range("E23") '  -0.12%  '
Set-TextValue $ws.Range("E24") '  +0.49%  '
Set-TextValue $ws.Range("E25") '  +2.36%  '
Set-TextValue $ws.Range("D26") '166.49'
Set-TextValue $ws.Range("E26") '  +0.58%  '
Set-TextValue $ws.Range("D27") '9.27'
Set-TextValue $ws.Range("E27") '  -0.31%  '
Set-TextValue $ws.Range("E28") '  +4.14%  '
Set-TextValue $ws.Range("E29") '  +1.14%  '
Set-TextValue $ws.Range("D30") '1.35'
Set-TextValue $ws.Range("E30") '  -0.61%  '
Set-TextValue $ws.Range("E31") '  +1.69%  '
Set-TextValue $ws.Range("E32") '  +0.16%  '
Set-TextValue $ws.Range("D33") '2.07'
Set-TextValue $ws.Range("E33") '  +12.06%  '
Set-TextValue $ws.Range("E34") '  +2.66%  '
Set-TextValue $ws.Range("E35") '  -0.79%  '
Set-TextValue $ws.Range("E36") '  -1.22%  '
Set-TextValue $ws.Range("E37") '  +8.59%  '
Set-TextValue $ws.Range("E38") '  +4.57%  '
Set-TextValue $ws.Range("E39") '  +0.06%  '
Set-TextValue $ws.Range("D40") '0.0218'
Set-TextValue $ws.Range("E40") '  +0.68%  '
Set-TextValue $ws.Range("E41") '  +1.52%  '
Set-TextValue $ws.Range("D42") '1.484.74'
Set-TextValue $ws.Range("E42") '  +0.32%  '
Set-TextValue $ws.Range("E43") '  +2.94%  '
Set-TextValue $ws.Range("E44") '  +1.38%  '
Set-TextValue $ws.Range("D45") '16.63'
Set-TextValue $ws.Range("E45") '  +1.46%  '
Set-TextValue $ws.Range("D46") '4.23'
Set-TextValue $ws.Range("E46") '  +16.65%  '
Set-TextValue $ws.Range("E47") '  -0.18%  '
Set-TextValue $ws.Range("E48") '  -0.22%  '
Set-TextValue $ws.Range("E49") '  +1.55%  '
Set-TextValue $ws.Range("D50") '6.99'
Set-TextValue $ws.Range("E50") '  -3.52%  '
Set-TextValue $ws.Range("D51") '2.244.31'
Set-TextValue $ws.Range("E51") '  +1.40%  '
